$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13: new entry for Secretaria General (Coprecovid announcement) ---
# Seed formatting for the new rows by copying the previous "last" data row (12),
# then overwrite with the real values so wrap/border/fill match the table body.
$ws.Range("A12:K12").Copy()
$ws.Range("A13:K13").PasteSpecial(-4122)
$ws.Range("A12:K12").Copy()
$ws.Range("A14:K14").PasteSpecial(-4122)

$ws.Range("A13").Value = "Secretaría General de la Presidencia de la República de Guatemala"
$ws.Range("B13").Formula = "=+ROW()-1"
$ws.Range("C13").Value = "General"
$ws.Range("D13").Value = "La Secretaría General de la Presidencia de la República (SGP) es el órgano responsable del apoyo jurídico y administrativo con carácter inmediato y constante del Presidente de la República. Su función es tramitar los asuntos de Gobierno del Despacho del Presidente."

$ws.Range("E13").Value = "https://legal.dca.gob.gt/"
$ws.Hyperlinks.Add($ws.Range("E13"), "https://legal.dca.gob.gt/GestionDocumento/VisualizarDocumento?verDocumentoPrevia=True&versionImpresa=False&doc=60004")
$ws.Range("E12").Copy()
$ws.Range("E13").PasteSpecial(-4122)

$ws.Range("F13").Value = "nuevas Disposiciones presidenciales en caso de calamidad pública `ny órdenes para el estricto cumplimiento, 24 mayo 2020`n"

$ws.Range("G13").Value = "https://sgp.gob.gt/"
$ws.Hyperlinks.Add($ws.Range("G13"), "https://sgp.gob.gt/")
$ws.Range("G12").Copy()
$ws.Range("G13").PasteSpecial(-4122)

$ws.Range("I13").Value = "25/5/2020"
$ws.Range("J13").Value = "Guatemala"

$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: new entry for Coprecovid creation agreement ---
$ws.Range("A14").Value = "Secretaría General de la Presidencia de la República de Guatemala"
$ws.Range("B14").Formula = "=+ROW()-1"
$ws.Range("C14").Value = "General"
$ws.Range("D14").Value = "La Secretaría General de la Presidencia de la República (SGP) es el órgano responsable del apoyo jurídico y administrativo con carácter inmediato y constante del Presidente de la República. Su función es tramitar los asuntos de Gobierno del Despacho del Presidente."

$ws.Range("E14").Value = "https://legal.dca.gob.gt/GestionDocumento/VisualizarDocumento?verDocumentoPrevia=True&versionImpresa=False&doc=60004"
$ws.Hyperlinks.Add($ws.Range("E14"), "https://legal.dca.gob.gt/GestionDocumento/VisualizarDocumento?verDocumentoPrevia=True&versionImpresa=False&doc=60004")
$ws.Range("E12").Copy()
$ws.Range("E14").PasteSpecial(-4122)

$ws.Range("F14").Value = "Acuerdo Gubernativo para la creación de la Comisión Presidencial de Atención a la Emergencia de Covid 19 ó Coprecovid"

$ws.Range("G14").Value = "https://sgp.gob.gt/"
$ws.Hyperlinks.Add($ws.Range("G14"), "https://sgp.gob.gt/")
$ws.Range("G12").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("I14").Value = "24/5/2020"
$ws.Range("J14").Value = "Guatemala"

$ws.Rows.Item(14).RowHeight = 91.5

# --- Extend the table + validation to cover the new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:K14"))

$dv = $ws.Range("C2:C12").Validation
$ws.Range("C2:C14").Validation.Delete()
$ws.Range("C2:C14").Validation.Add(3, 1, 1, "=Categoria")

# --- View state: scrolled/selected like the saved session ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("J15").Select()
